$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 (shifts existing rows 23:153 down to 24:154),
# matching the new weekly data point being prepended to the series.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new week's record.
$ws.Range("A23").Value = 8
$ws.Range("B23").Value = "Terminal La Palmera de La Serena"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44749
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = 100112040
$ws.Range("G23").Value = "Cilantro"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 2600
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = 1750
$ws.Range("N23").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O23").Value = "Provincia del Elquí"
$ws.Range("P23").Value = 1167
$ws.Range("Q23").Value = 1.5
$ws.Range("R23").Value = "Hortaliza"
